# Insert a new weekly price-report row for "Femacal de La Calera" /
# Albahaca right before the current row 213. This pushes the existing
# rows 213-249 down to 214-250 (data unchanged) and adds a brand new
# observation dated 2023-03-10 (serial 44995) as the new row 213.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 213:249 down to 214:250, leaving a blank row 213 behind
# (inherits formatting from the row it was split from, same as Excel's
# native Insert behaviour).
$ws.Rows.Item(213).EntireRow.Insert()

# Populate the new row 213 with the new observation.
$ws.Cells.Item(213, 1).Value  = 3
$ws.Cells.Item(213, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(213, 3).Value  = "Coquimbo"
$ws.Cells.Item(213, 4).Value  = 44995
$ws.Cells.Item(213, 5).Value  = 5
$ws.Cells.Item(213, 6).Value  = 100112052
$ws.Cells.Item(213, 7).Value  = "Albahaca"
$ws.Cells.Item(213, 8).Value  = "Sin especificar"
$ws.Cells.Item(213, 9).Value  = "Primera"
$ws.Cells.Item(213, 10).Value = 120
$ws.Cells.Item(213, 11).Value = 4500
$ws.Cells.Item(213, 12).Value = 5000
$ws.Cells.Item(213, 13).Value = 4729
$ws.Cells.Item(213, 14).Value = "$/docena de matas"
$ws.Cells.Item(213, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(213, 16).Value = 788
$ws.Cells.Item(213, 17).Value = 6
$ws.Cells.Item(213, 18).Value = "Hortaliza"
